$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "Save" - copy formatting from G1 (bold/border/center header style)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data column values (0/1 flags) for rows 2-7
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
